$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue "D2" "43.062.02"
Set-TextValue "E2" "  +0.17%  "
Set-TextValue "D3" "2.306.94"
Set-TextValue "E3" "  +0.16%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "300.06"
Set-TextValue "E5" "  -0.47%  "
Set-TextValue "D6" "97.94"
Set-TextValue "E6" "  -0.63%  "
Set-TextValue "E7" "  -1.47%  "
Set-TextValue "E8" "  +0.00%  "
Set-TextValue "E9" "  -2.55%  "
Set-TextValue "D10" "36.01"
Set-TextValue "E10" "  +0.91%  "
Set-TextValue "D11" "0.0792"
Set-TextValue "E11" "  +0.31%  "
Set-TextValue "E12" "  +1.37%  "
Set-TextValue "D13" "0.119"
Set-TextValue "E13" "  +1.93%  "
Set-TextValue "E14" "  -1.05%  "
Set-TextValue "D15" "2.666.90"
Set-TextValue "E15" "  +0.17%  "
Set-TextValue "D16" "2.308.62"
Set-TextValue "E16" "  +0.00%  "
Set-TextValue "E17" "  -0.93%  "
Set-TextValue "D18" "42.984.85"
Set-TextValue "E18" "  +0.19%  "
Set-TextValue "D19" "12.76"
Set-TextValue "E19" "  -5.18%  "
Set-TextValue "D21" "6.05"
Set-TextValue "E21" "  -1.51%  "
Set-TextValue "D22" "68.14"
Set-TextValue "E22" "  -0.02%  "
Set-TextValue "D23" "240.54"
Set-TextValue "E23" "  +0.65%  "
Set-TextValue "E24" "  -0.49%  "
Set-TextValue "E25" "  +0.11%  "
Set-TextValue "E26" "  +0.05%  "
Set-TextValue "E27" "  +0.11%  "
Set-TextValue "D28" "25.47"
Set-TextValue "E28" "  +2.72%  "
Set-TextValue "D29" "165.73"
Set-TextValue "E29" "  -0.81%  "
Set-TextValue "D30" "9.09"
Set-TextValue "E30" "  -0.43%  "
Set-TextValue "E31" "  -0.68%  "
Set-TextValue "D32" "33.19"
Set-TextValue "E32" "  -0.38%  "
Set-TextValue "D33" "4.98"
Set-TextValue "E33" "  +2.88%  "
Set-TextValue "E34" "  +0.04%  "
Set-TextValue "E35" "  -3.72%  "
Set-TextValue "D36" "17.07"
Set-TextValue "E36" "  -6.02%  "
Set-TextValue "E37" "  -0.92%  "
Set-TextValue "D38" "0.0687"
Set-TextValue "E38" "  -0.32%  "
Set-TextValue "E39" "  -0.66%  "
Set-TextValue "E40" "  -1.32%  "
Set-TextValue "E41" "  -0.26%  "
Set-TextValue "E42" "  -1.11%  "
Set-TextValue "D43" "2.015.33"
Set-TextValue "E43" "  +0.73%  "
Set-TextValue "E44" "  -2.11%  "
Set-TextValue "E45" "  +1.84%  "
Set-TextValue "D46" "10.12"
Set-TextValue "E46" "  -0.31%  "
Set-TextValue "D47" "17.43"
Set-TextValue "E47" "  -0.19%  "
Set-TextValue "E48" "  -0.81%  "
Set-TextValue "D49" "2.90"
Set-TextValue "E49" "  -3.84%  "
Set-TextValue "D50" "53.85"
Set-TextValue "E50" "  -1.39%  "
Set-TextValue "D51" "2.533.84"
Set-TextValue "E51" "  +0.17%  "
